# Report time for JC (Johan Can), NH (Noah Hellman), EH (Emir Hadzisalihovic)
# & YH (Yousef Hashem) in the timesheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- HELA GRUPPEN - TID block (hours worked) ---
# Johan Can, Maandag (week 41) hours -> 3.5 (feeds C6/M6/C11 totals via SUM formulas)
$ws.Range("C43").Value = 3.5

# --- HELA GRUPPEN - SYSSELSATTNING (AKTIVITETSNUMMER) block (activity numbers) ---
# Emir Hadzisalihovic, Maandag (week 41) activity number -> 2 (feeds C8/M8/C24 totals)
$ws.Range("Q17").Value = 2

# Yousef Hashem, Maandag (week 41) activity number -> 2 (feeds C9/M9/C37 totals)
$ws.Range("Q30").Value = 2

# Noah Hellman, Maandag (week 41) activity number -> 3 (feeds C10/M10/C50 totals)
$ws.Range("Q43").Value = 3

# Per-person weekly activity-number summary cells (week 41 column, "AE")
# Johan Can -> 3
$ws.Range("AE6").Value = 3
# Emir Hadzisalihovic -> "1, 2"
$ws.Range("AE8").Value = "1, 2"
# Yousef Hashem -> 6
$ws.Range("AE9").Value = 6
# Noah Hellman -> "2, 3"
$ws.Range("AE10").Value = "2, 3"

# Restore the cursor position recorded in the saved workbook.
$ws.Range("AE9").Select()
